$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F across rows 2-25
$valsBF = @(
    @(1.488102091308974,0.0429784381949645,0.137095416425641,0.07377937550791813,2.443551437920206),
    @(1.437926680809284,0.03891360204063687,0.137507215632855,0.0732840471825007,2.4171296269717),
    @(1.408023314977584,0.03639846191256879,0.1377667964018716,0.07301236143795009,2.401974831123155),
    @(1.396064652243638,0.03536859521847191,0.1378742713089851,0.07290980921766987,2.396067340362009),
    @(1.394092640746123,0.03519728704418412,0.1378922197356749,0.07289327354374997,2.395102593200022),
    @(1.407861116668016,0.03638459280985273,0.1377682389890431,0.0730109453315535,2.401894075227389),
    @(1.470613698455253,0.04158084816428698,0.1372360072340371,0.07360185295801003,2.434219202140582),
    @(1.600868671506589,0.05162101994585555,0.13624570200143,0.0750181162636494,2.506110257756092),
    @(1.700995526296424,0.05891186340332411,0.1355505948541254,0.07621589625807346,2.564154861022246),
    @(1.747517181800561,0.06221129379763113,0.1352414066785244,0.0767950204977943,2.591705720176549),
    @(1.765274246038643,0.06345831406019897,0.1351253347963599,0.07701924702137219,2.602303973589926),
    @(1.761443698059907,0.06318985162317858,0.1351502879620234,0.07697073679732114,2.60001408785871),
    @(1.748975251877823,0.06231393463566803,0.1352318371039001,0.07681336904942881,2.592574328078513),
    @(1.741356253126071,0.06177709948372012,0.1352819200158368,0.07671761810791722,2.588038808750127),
    @(1.697974843951954,0.0586958977258405,0.1355709424092826,0.07617873827331678,2.562377455234994),
    @(1.671611317661643,0.05680132354439138,0.1357500472105659,0.07585692457054094,2.546928962712968),
    @(1.656539318896705,0.05570999025265166,0.1358537241638729,0.07567504942661074,2.538151240075862),
    @(1.65145192195871,0.05534020237037396,0.1358889408132065,0.07561402328852651,2.535197759720901),
    @(1.674408276914733,0.05700317178715864,0.1357309128192874,0.07589084856096306,2.548562316091321),
    @(1.752633725649162,0.06257127730809486,0.1352078566939259,0.07685945811509498,2.594755074238748),
    @(1.804576655317305,0.06619638701275221,0.1348719028512058,0.07752120574656729,2.625908701538748),
    @(1.776778765894676,0.06426285183016489,0.1350506679081338,0.07716539201417305,2.609193039543925),
    @(1.673143507698853,0.05691192279809343,0.1357395612742351,0.07587550174506319,2.547823553671037),
    @(1.564856545662906,0.0489203363820252,0.1365078984170989,0.07460738526527422,2.485747211794575)
)

# New values for columns K,L across rows 2-25
$valsKL = @(
    @(1.006837424343303,0.233212589726449),
    @(0.9578911728307276,0.2262685050182256),
    @(0.9284730697622763,0.2221438623066518),
    @(0.9166442831922552,0.2204979433146974),
    @(0.9146897347119136,0.2202267462748182),
    @(0.9283128979756157,0.2221215236096441),
    @(0.9898287004303086,0.2307893697286829),
    @(1.115524746964752,0.2488942489684973),
    @(1.211005997462451,0.2628781551621557),
    @(1.255133650961909,0.2693896310846071),
    @(1.271943897005229,0.2718770480893795),
    @(1.268319052308442,0.2713403736908475),
    @(1.256514632061823,0.2695938378007838),
    @(1.249297128114051,0.2685268574545887),
    @(1.208136129247521,0.2624556421501723),
    @(1.183062933391824,0.2587696516748963),
    @(1.168706672391551,0.2566636935534916),
    @(1.163857070117587,0.2559530762532063),
    @(1.185725269026506,0.2591605691204393),
    @(1.259979157204526,0.2701062492168091),
    @(1.309091813167782,0.2773861845584946),
    @(1.282825929472892,0.2734891645520321),
    @(1.184521444425741,0.2589837943374533),
    @(1.080973873212145,0.2438771174270471)
)

# New values for column N across rows 2-25
$valsN = @(
    3.257016017891075,
    3.258522478605315,
    3.260171859267004,
    3.261025872337072,
    3.261178657333303,
    3.260182640766814,
    3.257384921364263,
    3.257660754437737,
    3.261398319846393,
    3.263871213444673,
    3.264919118760105,
    3.264688469847215,
    3.263955189016741,
    3.26352056154731,
    3.261252289515852,
    3.260058933434038,
    3.259445250397661,
    3.2592499437024,
    3.260178441059836,
    3.264167543040486,
    3.26742458282277,
    3.265626647997152,
    3.26012418617762,
    3.256966985954421
)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = $valsBF[$i][0]
    $ws.Cells.Item($row, 3).Value = $valsBF[$i][1]
    $ws.Cells.Item($row, 4).Value = $valsBF[$i][2]
    $ws.Cells.Item($row, 5).Value = $valsBF[$i][3]
    $ws.Cells.Item($row, 6).Value = $valsBF[$i][4]

    $ws.Cells.Item($row, 11).Value = $valsKL[$i][0]
    $ws.Cells.Item($row, 12).Value = $valsKL[$i][1]

    $ws.Cells.Item($row, 14).Value = $valsN[$i]
}
